$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a typo: " 은 즉" -> "는 즉" for the "rightarrow" row
$ws.Range("B68").Value = "는 즉"

# Remove the "(" / ")" rows (괄호열고/괄호닫고) - rows 73 and 74
$ws.Rows.Item(73).Delete()
$ws.Rows.Item(73).Delete()

# Append two new rows at the bottom of the table for the absolute-value bars
$ws.Range("A96").Value = "left|"
$ws.Range("B96").Value = " 절댓값"
$ws.Range("A97").Value = "right|"
$ws.Range("B97").Value = " "

# Update the view/selection to match the final state
$ws.Range("D79").Select()
